$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.359664797782898
$ws.Range("B1").Value = 1.68866229057312
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.820631742477417
$ws.Range("E1").Value = 0.7754678130149841
